$d = $word.ActiveDocument

# Helper: find $searchText, and replace the range (plus one character of
# context on either side) with $before + $replaceText + $after. Including
# a little context on both sides of the match ensures any proofing marks
# (w:proofErr) that sit at the edges of the matched run are cleared away
# along with the run, instead of being left behind as stray/empty markers.
function Replace-WithContext($searchText, $replaceText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) {
        return $false
    }
    $s = $rng.Start
    $e = $rng.End
    $before = $d.Range($s - 1, $s)
    $after = $d.Range($e, $e + 1)
    $wide = $d.Range($s - 1, $e + 1)
    $wide.Text = $before.Text + $replaceText + $after.Text
    return $true
}

# ------------------------------------------------------------------
# 1. Summary paragraph edits
# ------------------------------------------------------------------

# "Strives for thoroughness" -> "Aims for thoroughness"
$rng = $d.Content
$rng.Find.Execute("Strives for thoroughness", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "Aims for thoroughness", 2)

# "highly-consistent" -> "highly consistent" (also clears the gramStart/gramEnd
# proofErr markers that wrapped the old text)
Replace-WithContext "highly-consistent" "highly consistent"

# "design and develop software systems" -> "design and develop highly scalable software systems"
$rng = $d.Content
$rng.Find.Execute("design and develop software systems", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "design and develop highly scalable software systems", 2)

# ------------------------------------------------------------------
# 2. Move the _GoBack bookmark out of the empty paragraph (after the
#    "Certificate in Computer Programming" bullet) and into the summary
#    paragraph, right before "highly scalable".
# ------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$rng = $d.Content
$rng.Find.Execute("design and develop highly scalable", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)
$insertPos = $rng.Start + ("design and develop ").Length
$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 3. "Contributed to widely-known" paragraph - collapse the separate
#    runs (and remove the gramStart/gramEnd proofErr markers around
#    "widely-known") into a single piece of text. The wording itself is
#    unchanged, so round-trip through a placeholder to force the engine
#    to actually rebuild the run (and drop the now orphaned proofErr
#    markers) rather than treat it as a no-op.
# ------------------------------------------------------------------

Replace-WithContext "widely-known" "TEMP_PLACEHOLDER_WIDELY_KNOWN"
$rng = $d.Content
$rng.Find.Execute("TEMP_PLACEHOLDER_WIDELY_KNOWN", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "widely-known", 2)
